# Add a new "Trailer1Axle_f" droplink hardpoint sheet to the workbook,
# modeled on the existing "Bus_Makulu_r" sheet (same layout/styles),
# with updated hardpoint values, and make it the active sheet.

$wb = $excel.ActiveWorkbook

# Use the last existing sheet ("Bus_Makulu_r") as the template - it has
# identical layout/conditional-formatting/styles to the new sheet.
$template = $wb.Worksheets.Item("Bus_Makulu_r")

# Copy it to the end of the workbook (after the last sheet).
$template.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null

# The copy becomes the sheet right after the template, i.e. the new last sheet.
$new = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Name = "Trailer1Axle_f"

# New instance label for this hardpoint set.
$new.Range("H3").Value = "Droplink_Trailer1Axle_f"

# Updated sOutboard hardpoint (row 5).
$new.Range("F5").Value = 0.05
$new.Range("G5").Value = 0.6
$new.Range("H5").Value = 0.19

# Updated sInboard hardpoint (row 6) - x is a formula.
$new.Range("F6").Formula = "=0.3-0.15"
$new.Range("G6").Value = 0.58
$new.Range("H6").Value = 0.2

# k and m (rows 7-8) keep the same values as the template (50, 0.5).

# Make the new sheet the active one, with H7 selected (matches authored file).
$new.Activate() | Out-Null
$new.Range("D33").Select() | Out-Null
$new.Range("H7").Select() | Out-Null
